$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D sometimes hold numeric-looking text (e.g. "549.64").
# Excel auto-converts such strings to real numbers on assignment, which
# would corrupt the intended text formatting (and float precision). Force
# those cells to Text format before writing so the value stays a string,
# matching the original workbook layout.

$ws.Range("D2").Value = '57.691.05'
$ws.Range("E2").Value = '  -4.09%  '
$ws.Range("D3").Value = '2.935.21'
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '549.64'
$ws.Range("E5").Value = '  -4.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.54'
$ws.Range("E6").Value = '  +4.41%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.512'
$ws.Range("D9").Value = '2.928.73'
$ws.Range("E9").Value = '  -2.15%  '
$ws.Range("E10").Value = '  -3.97%  '
$ws.Range("E11").Value = '  -5.31%  '
$ws.Range("E12").Value = '  +1.30%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000221'
$ws.Range("E13").Value = '  +0.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.93'
$ws.Range("E14").Value = '  +1.37%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.121'
$ws.Range("E15").Value = '  +0.18%  '
$ws.Range("D16").Value = '3.421.56'
$ws.Range("E16").Value = '  -1.97%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.84'
$ws.Range("E17").Value = '  +6.58%  '
$ws.Range("D18").Value = '2.929.37'
$ws.Range("E18").Value = '  -2.14%  '
$ws.Range("D19").Value = '57.716.95'
$ws.Range("E19").Value = '  -3.99%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '417.78'
$ws.Range("E20").Value = '  -2.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.18'
$ws.Range("E21").Value = '  +0.67%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.685'
$ws.Range("E22").Value = '  +2.81%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.97'
$ws.Range("E23").Value = '  -0.96%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.03'
$ws.Range("E24").Value = '  +1.53%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '79.79'
$ws.Range("E25").Value = '  +0.63%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("E28").Value = '  -2.69%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.48'
$ws.Range("E29").Value = '  +3.30%  '
$ws.Range("E30").Value = '  +1.77%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '25.12'
$ws.Range("E31").Value = '  -0.34%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.99'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0968'
$ws.Range("E33").Value = '  +2.78%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.65'
$ws.Range("E34").Value = '  +1.11%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.937'
$ws.Range("E35").Value = '  +0.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.07'
$ws.Range("E36").Value = '  +1.06%  '
$ws.Range("E37").Value = '  -4.67%  '
$ws.Range("B38").Value = 'PEPE'
$ws.Range("C38").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D38").Value = '0.0₃0682'
$ws.Range("E38").Value = '  +2.39%  '
$ws.Range("B39").Value = 'Cosmos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.71'
$ws.Range("E39").Value = '  +3.23%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.55'
$ws.Range("E40").Value = '  +3.58%  '
$ws.Range("E41").Value = '  -0.20%  '
$ws.Range("E42").Value = '  +0.36%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0345'
$ws.Range("E43").Value = '  -2.82%  '
$ws.Range("D44").Value = '2.682.43'
$ws.Range("E44").Value = '  +0.21%  '
$ws.Range("E45").Value = '  +0.04%  '
$ws.Range("E46").Value = '  +1.71%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '122.02'
$ws.Range("E47").Value = '  +1.29%  '
$ws.Range("E48").Value = '  +1.54%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.97'
$ws.Range("E49").Value = '  -1.21%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.15'
$ws.Range("E50").Value = '  -1.61%  '
$ws.Range("E51").Value = '  -0.19%  '
